$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Monday block
$ws.Range("A2").Value = "Math1-1"
$ws.Range("B2").Value = "Math1-2"
$ws.Range("C2").Value = "Phys1-3"
$ws.Range("E2").Value = "Math1-5"
$ws.Range("A3").Value = "Russian1-1"
$ws.Range("B3").Value = "Russian1-2"
$ws.Range("C3").Value = "English1-3"
$ws.Range("D3").Value = "English1-4"
$ws.Range("D4").Value = "Math1-4"
$ws.Range("E4").Value = "Russian1-5"

# Tuesday block
$ws.Range("A7").Value = "English2-1"
$ws.Range("B7").Value = "Math2-2"
$ws.Range("C7").Value = "Math2-3"
$ws.Range("D7").Value = "English2-4"
$ws.Range("A8").Value = "Math2-1"
$ws.Range("B8").Value = "English2-2"
$ws.Range("C8").Value = "Phys2-3"
$ws.Range("D8").Value = "Phys2-4"

# Wednesday block
$ws.Range("A11").Value = "English3-1"
$ws.Range("B11").Value = "Russian3-2"
$ws.Range("D11").Value = "Phys3-4"
$ws.Range("A12").Value = "Phys3-1"
$ws.Range("B12").Value = "Phys3-2"
$ws.Range("C12").Value = "Phys3-3"
$ws.Range("D12").Value = "Litra3-4"
$ws.Range("A13").Value = "Russian3-1"

# Thursday block
$ws.Range("A16").Value = "Math4-1"
$ws.Range("A17").Value = "English4-1"
$ws.Range("B17").Value = "English4-2"

# Friday block
$ws.Range("B20").Value = "Russian5-2"
$ws.Range("D20").Value = "Russian5-4"
